$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "BloodPressure"
$ws.Range("C2").Value = 58
$ws.Range("C3").Value = 75
$ws.Range("C4").Value = 64
$ws.Range("C5").Value = 80
$ws.Range("C6").Value = 90
$ws.Range("C7").Value = 72
$ws.Range("C8").Value = 48
$ws.Range("C9").Value = 35
$ws.Range("C10").Value = 72
$ws.Range("C11").Value = 78
$ws.Range("C12").Value = 70
$ws.Range("C13").Value = 95
$ws.Range("C14").Value = 60
$ws.Range("C15").Value = 74
$ws.Range("C16").Value = 58
$ws.Range("C17").Value = 50
$ws.Range("C18").Value = 72
$ws.Range("C19").Value = 68
$ws.Range("C20").Value = 35
$ws.Range("C21").Value = 65
$ws.Range("C22").Value = 74
$ws.Range("C23").Value = 78
$ws.Range("C24").Value = 72
$ws.Range("C25").Value = 64
$ws.Range("C26").Value = 76
$ws.Range("C27").Value = 84
$ws.Range("C28").Value = 60
$ws.Range("C29").Value = 74
$ws.Range("C30").Value = 86
$ws.Range("C31").Value = 74
$ws.Range("C32").Value = 76
$ws.Range("D32").Value = 50.55
$ws.Range("C33").Value = 78
$ws.Range("F33").Value = 66.5
$ws.Range("C34").Value = 78
$ws.Range("C35").Value = 35
$ws.Range("C36").Value = 92
$ws.Range("C37").Value = 78
$ws.Range("C38").Value = 70
$ws.Range("C39").Value = 88
$ws.Range("C40").Value = 78
$ws.Range("C41").Value = 35
$ws.Range("C42").Value = 58
$ws.Range("C43").Value = 54
$ws.Range("C44").Value = 75
$ws.Range("C45").Value = 72
$ws.Range("C46").Value = 92
$ws.Range("C47").Value = 76
$ws.Range("C48").Value = 72
$ws.Range("E48").Value = 1.2
$ws.Range("C49").Value = 72
$ws.Range("C50").Value = 72
$ws.Range("C51").Value = 78
$ws.Range("C52").Value = 65
$ws.Range("C53").Value = 88
$ws.Range("C54").Value = 82
$ws.Range("C55").Value = 70
$ws.Range("C56").Value = 62
$ws.Range("C57").Value = 68
$ws.Range("C58").Value = 70
$ws.Range("C59").Value = 35
$ws.Range("D59").Value = 13.35
$ws.Range("C60").Value = 70
$ws.Range("C61").Value = 64
$ws.Range("C62").Value = 78
$ws.Range("C63").Value = 70
$ws.Range("C64").Value = 90
$ws.Range("C65").Value = 62
$ws.Range("C66").Value = 82
$ws.Range("C67").Value = 85
$ws.Range("C68").Value = 60
$ws.Range("C69").Value = 107
$ws.Range("C70").Value = 76
$ws.Range("C71").Value = 68
$ws.Range("C72").Value = 66
$ws.Range("C73").Value = 86
$ws.Range("C74").Value = 68
$ws.Range("C75").Value = 85
$ws.Range("C76").Value = 66
$ws.Range("C77").Value = 107
$ws.Range("C78").Value = 75
$ws.Range("C79").Value = 56
$ws.Range("C80").Value = 66
$ws.Range("C81").Value = 92
$ws.Range("C82").Value = 66
$ws.Range("C84").Value = 64
$ws.Range("C85").Value = 94
$ws.Range("C86").Value = 64
$ws.Range("C87").Value = 70
$ws.Range("C88").Value = 88
$ws.Range("C89").Value = 68
$ws.Range("C90").Value = 74
$ws.Range("C91").Value = 58
$ws.Range("C92").Value = 60
$ws.Range("C93").Value = 80
$ws.Range("C94").Value = 35
$ws.Range("D94").Value = 13.35
$ws.Range("C95").Value = 35
$ws.Range("C96").Value = 48
$ws.Range("C97").Value = 82
$ws.Range("C98").Value = 52
$ws.Range("C99").Value = 65
$ws.Range("C100").Value = 60
$ws.Range("C101").Value = 78
$ws.Range("C102").Value = 82
$ws.Range("C103").Value = 88
$ws.Range("C104").Value = 90
$ws.Range("C105").Value = 76
$ws.Range("C106").Value = 78
$ws.Range("C107").Value = 86
$ws.Range("C108").Value = 58
$ws.Range("C109").Value = 78
$ws.Range("C110").Value = 82
$ws.Range("E110").Value = 1.2
$ws.Range("C111").Value = 80
$ws.Range("C112").Value = 60
$ws.Range("C113").Value = 88
$ws.Range("C114").Value = 35
$ws.Range("D114").Value = 13.35
$ws.Range("C115").Value = 106
$ws.Range("C116").Value = 72
$ws.Range("C117").Value = 60
$ws.Range("C118").Value = 52
$ws.Range("C119").Value = 66
$ws.Range("C120").Value = 62
$ws.Range("C121").Value = 35
$ws.Range("C122").Value = 58
$ws.Range("C123").Value = 70
$ws.Range("C124").Value = 80
$ws.Range("C125").Value = 62
$ws.Range("E125").Value = 1.2
$ws.Range("C126").Value = 60
$ws.Range("C127").Value = 76
$ws.Range("C128").Value = 35
$ws.Range("D128").Value = 13.35
$ws.Range("C129").Value = 76
$ws.Range("E129").Value = 1.2
$ws.Range("C130").Value = 35
$ws.Range("D130").Value = 50.55
$ws.Range("C131").Value = 107
$ws.Range("C132").Value = 64
$ws.Range("C133").Value = 70
$ws.Range("C134").Value = 60
$ws.Range("C135").Value = 84
$ws.Range("C136").Value = 78
$ws.Range("E136").Value = 1.2
$ws.Range("C137").Value = 74
$ws.Range("C138").Value = 86
$ws.Range("C139").Value = 107
$ws.Range("C140").Value = 78
$ws.Range("C141").Value = 62
$ws.Range("C142").Value = 66
$ws.Range("C143").Value = 60
$ws.Range("C144").Value = 75
$ws.Range("C145").Value = 62
$ws.Range("C146").Value = 60
$ws.Range("C147").Value = 74
$ws.Range("C148").Value = 74
$ws.Range("E148").Value = 1.2
$ws.Range("C149").Value = 50
$ws.Range("C150").Value = 78
$ws.Range("C151").Value = 88
$ws.Range("E152").Value = 1.2
$ws.Range("C153").Value = 72
$ws.Range("C154").Value = 70
$ws.Range("C155").Value = 70
